# Daily attendance processing - 2025-12-19 10:31:08
# Normalises the "Recorded By" column (G) so that the automated/system
# recorder tag(s) are listed after the real user(s) instead of before.
#   "System, user@example.com"          -> "user@example.com, System"
#   "System, system, user@example.com"  -> "user@example.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ",\s*"

    # Find the leading run of "System"/"system" tokens.
    $leadCount = 0
    while (($leadCount -lt $parts.Length) -and ($parts[$leadCount].ToLower() -eq "system")) {
        $leadCount++
    }

    # Only reorder when there's a leading System run AND something after it.
    if (($leadCount -gt 0) -and ($leadCount -lt $parts.Length)) {
        $leading = $parts[0..($leadCount - 1)]
        $rest = $parts[$leadCount..($parts.Length - 1)]
        $newParts = $rest + $leading
        $cell.Value = [string]::Join(", ", $newParts)
    }
}
